$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: C5 gets new text and loses its yellow-fill style
$ws.Range("C5").Value = "@james_frain Loving James in this. He's making me giggle, which is plain wrong LOL"
$ws.Range("C5").ClearFormats()

# Row 11: C11 text updated
$ws.Range("C11").Value = "📷 Now if the wind changes… Hilarious photo found on Instagram."

# Row 7: C7 text updated
$ws.Range("C7").Value = "Well we girls didn’t give up like the boys did. @ Wombatz Retreat"

# Row 2: C2 text updated
$ws.Range("C2").Value = "Toys galore. A not so typical house in Boorowa. @ Boorowa, New South Wales"

# Row 14: C14 new text, loses red-fill style; D14 becomes a date value with the
# same format style used elsewhere in column D (copy format from D4), losing its
# red-fill style too.
$ws.Range("C14").Value = "@startrekcbs OMG Burnham. OMG Tyler. OMG Lorca. OMG L’rell. OMG Tilly. OMG Saru. OMG Georgiou. OMG Stamets. OMG Cul…"
$ws.Range("C14").ClearFormats()

$ws.Range("D4").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = 43122.497442129628

# Row 4: B4 numeric value changed
$ws.Range("B4").Value = 8489

# Update the active cell selection to D14
$ws.Range("D14").Select()
